$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vtab = [char]11

$cell = $t.Cell(1, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "80 x 73" + $vtab + "  7    3" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "0|    |"

$cell = $t.Cell(1, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "35 x 75" + $vtab + "  7    5" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "5|    |"

$cell = $t.Cell(1, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "80 x 60" + $vtab + "  6    0" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "0|    |"

$cell = $t.Cell(2, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "30 x 58" + $vtab + "  5    8" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "0|    |"

$cell = $t.Cell(2, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "53 x 19" + $vtab + "  1    9" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "3|    |"

$cell = $t.Cell(2, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "55 x 35" + $vtab + "  3    5" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "5|    |"

$cell = $t.Cell(3, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "99 x 86" + $vtab + "  8    6" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "9|    |"

$cell = $t.Cell(3, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "49 x 18" + $vtab + "  1    8" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "9|    |"

$cell = $t.Cell(3, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "94 x 34" + $vtab + "  3    4" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "4|    |"

$cell = $t.Cell(4, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "35 x 29" + $vtab + "  2    9" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "5|    |"

$cell = $t.Cell(4, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "43 x 29" + $vtab + "  2    9" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "3|    |"

$cell = $t.Cell(4, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "63 x 49" + $vtab + "  4    9" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "3|    |"

$cell = $t.Cell(5, 1)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "34 x 42" + $vtab + "  4    2" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "4|    |"

$cell = $t.Cell(5, 2)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "71 x 96" + $vtab + "  9    6" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "1|    |"

$cell = $t.Cell(5, 3)
$rng = $cell.Range
$rng.End = $rng.End - 1
$rng.Text = "57 x 23" + $vtab + "  2    3" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "7|    |"
